$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.575.03"
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").Value = "2.316.86"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "'106.17"
$ws.Range("E5").Value = "  +8.48%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'310.96"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +8.19%  "
$ws.Range("D10").Value = "'36.89"
$ws.Range("E10").Value = "  +5.13%  "
$ws.Range("D11").Value = "'53.14"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").Value = "'0.0819"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").Value = "'7.06"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").Value = "2.672.15"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").Value = "'15.24"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").Value = "2.314.38"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").Value = "'0.817"
$ws.Range("E18").Value = "  +3.63%  "
$ws.Range("D19").Value = "43.481.39"
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("D20").Value = "'12.25"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "0.0₃0933"
$ws.Range("E21").Value = "  +3.14%  "
$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = "  +4.19%  "
$ws.Range("D23").Value = "'68.41"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "'243.34"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("D25").Value = "'2.05"
$ws.Range("E25").Value = "  +3.77%  "
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'25.22"
$ws.Range("E28").Value = "  +6.83%  "
$ws.Range("E29").Value = "  +11.92%  "
$ws.Range("D30").Value = "'37.33"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "'165.70"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'18.39"
$ws.Range("E35").Value = "  +4.51%  "
$ws.Range("E36").Value = "  +6.61%  "
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("D38").Value = "'3.07"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").Value = "'4.56"
$ws.Range("E39").Value = "  +9.61%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.89"
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.107"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("E42").Value = "  +22.76%  "
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("D45").Value = "2.001.93"
$ws.Range("E45").Value = "  +2.81%  "
$ws.Range("D46").Value = "'3.16"
$ws.Range("E46").Value = "  +6.86%  "
$ws.Range("D47").Value = "'19.14"
$ws.Range("E47").Value = "  +2.46%  "
$ws.Range("D48").Value = "'10.04"
$ws.Range("E48").Value = "  +2.51%  "
$ws.Range("D49").Value = "'57.94"
$ws.Range("E49").Value = "  +6.03%  "
$ws.Range("E50").Value = "  +9.45%  "
$ws.Range("D51").Value = "'2.91"
$ws.Range("E51").Value = "  +0.93%  "
